$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Persons")

# ------------------------------------------------------------------
# Copy the header formatting (fill/font/alignment - style indexes 3
# and 4) from the old header row (row 1) onto the new header row
# (row 4, shifted one column to the right) before the old row is
# cleared, so Excel re-uses the existing style records instead of
# minting new ones.
# ------------------------------------------------------------------
$ws.Range("A1:C1").Copy()
$ws.Range("B4:D4").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("D1:G1").Copy()
$ws.Range("E4:H4").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

# --- drop the old header row and the now-superseded placeholder rows ---
$ws.Range("A1:G1").Clear()
$ws.Range("D3:G3").Clear()
$ws.Range("D5:G5").Clear()

# --- new header row text (row 4, columns B:H) ---
$ws.Range("B4").Value = "Name"
$ws.Range("C4").Value = "Surname"
$ws.Range("D4").Value = "MarriedYear"
$ws.Range("E4").Value = "BirthDate"
$ws.Range("F4").Value = "OwnCar.Name"
$ws.Range("G4").Value = "OwnCar.Targa"
$ws.Range("H4").Value = "OwnCar.BuildYear"

# --- new data row (row 5, columns B:H) ---
$ws.Range("B5").Value = "pippo"
$ws.Range("C5").Value = "pluto"

$ws.Range("D5").HorizontalAlignment = -4152   # xlRight
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "2000"

$ws.Range("E5").NumberFormat = "mm-dd-yy"
$ws.Range("E5").Value = [datetime]"1980-10-10"

$ws.Range("F5").Value = "Audi"
$ws.Range("G5").Value = "DQ789AQ"
$ws.Range("H5").Value = 2009

# --- new columns H and I get the same width as column G ---
$ws.Range("H1:I1").EntireColumn.ColumnWidth = $ws.Range("G1").EntireColumn.ColumnWidth

# --- selection / active sheet ---
$ws.Range("B5:H5").Select()
$ws.Activate()
